$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.015833333333333
$ws.Range("H2").Value = 9.047499999999999
$ws.Range("I2").Value = 0.05376901095572644
$ws.Range("J2").Value = 0.05376901095572643
$ws.Range("M2").Value = 12.67919733333333
$ws.Range("N2").Value = 38.037592
$ws.Range("O2").Value = 0.9871416146107245
$ws.Range("P2").Value = 0.9871416146107247
$ws.Range("Q2").Value = 38.23834595777778
$ws.Range("R2").Value = 344.14511362
$ws.Range("S2").Value = 0.05307762829085753
$ws.Range("T2").Value = 0.05307762829085754

$ws.Range("G3").Value = 3.015833333333333
$ws.Range("H3").Value = 9.047499999999999
$ws.Range("I3").Value = 0.05376901095572644
$ws.Range("J3").Value = 0.05376901095572643
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.01285838538927542
$ws.Range("P3").Value = 0.01285838538927542
$ws.Range("Q3").Value = 0.4980879963888889
$ws.Range("R3").Value = 4.4827919675
$ws.Range("S3").Value = 0.0006913826648689027
$ws.Range("T3").Value = 0.0006913826648689026

$ws.Range("I4").Value = 0.03577730856453667
$ws.Range("J4").Value = 0.03577730856453666
$ws.Range("M4").Value = 12.67919733333333
$ws.Range("N4").Value = 38.037592
$ws.Range("O4").Value = 0.9871416146107245
$ws.Range("P4").Value = 0.9871416146107247
$ws.Range("Q4").Value = 25.44337487359378
$ws.Range("R4").Value = 228.990373862344
$ws.Range("S4").Value = 0.03531727014282283
$ws.Range("T4").Value = 0.03531727014282283

$ws.Range("I5").Value = 0.03577730856453667
$ws.Range("J5").Value = 0.03577730856453666
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.01285838538927542
$ws.Range("P5").Value = 0.01285838538927542
$ws.Range("Q5").Value = 0.331422275067889
$ws.Range("S5").Value = 0.0004600384217138366
$ws.Range("T5").Value = 0.0004600384217138364

$ws.Range("G6").Value = 1.732509666666666
$ws.Range("H6").Value = 5.197528999999999
$ws.Range("I6").Value = 0.03088875310789786
$ws.Range("J6").Value = 0.03088875310789786
$ws.Range("M6").Value = 12.67919733333333
$ws.Range("N6").Value = 38.037592
$ws.Range("O6").Value = 0.9871416146107245
$ws.Range("P6").Value = 0.9871416146107247
$ws.Range("Q6").Value = 21.96683194557422
$ws.Range("R6").Value = 197.701487510168
$ws.Range("S6").Value = 0.03049157361624232
$ws.Range("T6").Value = 0.03049157361624233

$ws.Range("G7").Value = 1.732509666666666
$ws.Range("H7").Value = 5.197528999999999
$ws.Range("I7").Value = 0.03088875310789786
$ws.Range("J7").Value = 0.03088875310789786
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.01285838538927542
$ws.Range("P7").Value = 0.01285838538927542
$ws.Range("Q7").Value = 0.2861372540241111
$ws.Range("R7").Value = 2.575235286217
$ws.Range("S7").Value = 0.0003971794916555294
$ws.Range("T7").Value = 0.0003971794916555294

$ws.Range("G8").Value = 49.33364366666667
$ws.Range("H8").Value = 148.000931
$ws.Range("I8").Value = 0.879564927371839
$ws.Range("J8").Value = 0.879564927371839
$ws.Range("M8").Value = 12.67919733333333
$ws.Range("N8").Value = 38.037592
$ws.Range("O8").Value = 0.9871416146107245
$ws.Range("P8").Value = 0.9871416146107247
$ws.Range("Q8").Value = 625.5110032220169
$ws.Range("R8").Value = 5629.599028998153
$ws.Range("S8").Value = 0.8682551425608017
$ws.Range("T8").Value = 0.8682551425608019

$ws.Range("G9").Value = 49.33364366666667
$ws.Range("H9").Value = 148.000931
$ws.Range("I9").Value = 0.879564927371839
$ws.Range("J9").Value = 0.879564927371839
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.01285838538927542
$ws.Range("P9").Value = 0.01285838538927542
$ws.Range("Q9").Value = 8.147829476151445
$ws.Range("R9").Value = 73.330465285363
$ws.Range("S9").Value = 0.01130978481103715
$ws.Range("T9").Value = 0.01130978481103715
